# Add an "Email" column (C) to the "Users" worksheet, with a sample
# admin email address, matching the bold header style used by the
# existing Username/Password headers.

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Users")

# Header cell C1 - copy formatting from A1 (bold header style) onto C1,
# then set its text to "Email".
$ws.Range("A1").Copy()
$ws.Range("C1").PasteSpecial(-4122)  # xlPasteFormats
$ws.Range("C1").Value = "Email"

# Data cell C2 - sample admin email
$ws.Range("C2").Value = "admin@example.com"
